# APS1 - TransCal : versao atualizada do codigo
# Expand the truss model data (nodes, members, loads, supports) across all
# four worksheets and tidy up a couple of derived "count" cells.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Nos" (node coordinates)
# ---------------------------------------------------------------------------
$wsNos = $wb.Worksheets.Item("Nos")

$nosA = @(0,1,2,3,4,1,2,3,2)
$nosB = @(0,0,0,0,0,1,1,1,2)
for ($i = 0; $i -lt $nosA.Length; $i++) {
    $r = $i + 2
    $wsNos.Cells.Item($r, 1).Value = $nosA[$i]
    $wsNos.Cells.Item($r, 2).Value = $nosB[$i]
}

# D2 used to hold =COUNT(A2:A1048576); it's now a plain value.
$wsNos.Range("D2").Value = 9

# D9 picks up the (new) underlined-font style with no value.
$wsNos.Range("D9").Font.Underline = $true

# ---------------------------------------------------------------------------
# Sheet "Incidencia" (member connectivity)
# ---------------------------------------------------------------------------
$wsInc = $wb.Worksheets.Item("Incidencia")

$incA = @(1,2,3,4,1,2,2,3,4,6,7,4,5,6,7,8)
$incB = @(2,3,4,5,6,6,7,7,7,7,8,8,8,9,9,9)
for ($i = 0; $i -lt $incA.Length; $i++) {
    $r = $i + 2
    $wsInc.Cells.Item($r, 1).Value = $incA[$i]
    $wsInc.Cells.Item($r, 2).Value = $incB[$i]
    $wsInc.Cells.Item($r, 3).Value = 210000000000
}

$wsInc.Range("D2:D17").Formula = "=0.15*0.15"

# F2 used to hold =COUNT(A2:A1048576); it's now a plain value.
$wsInc.Range("F2").Value = 16

# F18 picks up the (new) underlined-font style with no value.
$wsInc.Range("F18").Font.Underline = $true

# ---------------------------------------------------------------------------
# Sheet "Carregamento" (loads)
# ---------------------------------------------------------------------------
$wsCar = $wb.Worksheets.Item("Carregamento")

$wsCar.Range("A2").Value = 6
$wsCar.Range("B2").Value = 1
$wsCar.Range("C2").Value = 2000

$wsCar.Range("A3").Value = 9
$wsCar.Range("B3").Value = 1
$wsCar.Range("C3").Value = 2000

$wsCar.Range("A4").Value = 9
$wsCar.Range("B4").Value = 2
$wsCar.Range("C4").Value = 10000

# ---------------------------------------------------------------------------
# Sheet "Restricao" (supports)
# ---------------------------------------------------------------------------
$wsRes = $wb.Worksheets.Item("Restricao")

$resA = @(1,1,2,3,4,5,5)
$resB = @(1,2,2,2,2,1,2)
for ($i = 0; $i -lt $resA.Length; $i++) {
    $r = $i + 2
    $wsRes.Cells.Item($r, 1).Value = $resA[$i]
    $wsRes.Cells.Item($r, 2).Value = $resB[$i]
}

# ---------------------------------------------------------------------------
# Active sheet / selections (matches the saved view state in the workbook)
# ---------------------------------------------------------------------------
$wsNos.Range("C22").Select()
$wsCar.Range("F10").Select()
$wsRes.Range("C13").Select()

$wsInc.Activate()
$wsInc.Range("C22").Select()
